# Features_Status.xlsx -- BoardManager.c row marked Done, status notes updated
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 5 (Console): fill in the "Left to do" cell with "Done" ---
$ws.Range("H5").Value = "Done"

# --- Row 6 (BoardManager.c): mark as finished ---
# Assignee is now both Tomer and Or
$ws.Range("G6").Value = "Tomer / Or"
# Replace the old "left to do" note with a short status update
$ws.Range("H6").Value = "All functions updated. No tests done."

# Re-color the whole BoardManager.c row from yellow (in-progress) to green (done)
$ws.Range("E6:H6").Interior.Color = $ws.Range("H5").Interior.Color
$ws.Range("H6").WrapText = $true

# Row no longer needs the tall wrapped height it had while the note was long
$ws.Rows.Item(6).AutoFit()

# --- View: scroll the sheet one column to the right (topLeftCell E1 -> F1) ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1

# Keep the original selection
$ws.Range("H7").Select()
